$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("C4").Value = -13.3512
$ws.Range("B9").Value = 5.885000000000008
$ws.Range("C9").Value = -12.31590000000001
$ws.Range("D9").Value = -7.644200000000005
$ws.Range("C11").Value = -13.6794
$ws.Range("B13").Value = 6.342500000000001
$ws.Range("B16").Value = 4.939299999999999
$ws.Range("C16").Value = -13.7762
$ws.Range("B18").Value = 6.415899999999996
$ws.Range("B20").Value = 8.263100000000003
$ws.Range("D22").Value = -8.112200000000007
$ws.Range("C23").Value = -12.9652
$ws.Range("C24").Value = -12.5203
$ws.Range("B26").Value = 5.358300000000002
$ws.Range("C26").Value = -12.96450000000001
$ws.Range("B27").Value = 6.170700000000007
$ws.Range("D27").Value = -7.858600000000003
$ws.Range("B29").Value = 5.211900000000001
$ws.Range("D29").Value = -7.437099999999988
$ws.Range("D32").Value = -6.811399999999993
$ws.Range("C34").Value = -11.87620000000001
$ws.Range("B35").Value = 8.042700000000007
$ws.Range("C35").Value = -13.09980000000002
$ws.Range("B36").Value = 8.779000000000009
$ws.Range("D37").Value = -7.607800000000002
$ws.Range("D38").Value = -7.581600000000005
$ws.Range("D39").Value = -8.02859999999999
$ws.Range("D41").Value = -8.221800000000004
$ws.Range("C44").Value = -13.2017
$ws.Range("B45").Value = 4.742100000000005
$ws.Range("D45").Value = -7.480999999999995
$ws.Range("C48").Value = -12.7066
$ws.Range("D48").Value = -8.156500000000001
$ws.Range("C49").Value = -13.7304
$ws.Range("D51").Value = -8.748000000000003
$ws.Range("C52").Value = -11.1175
$ws.Range("B55").Value = 6.894699999999996
$ws.Range("D56").Value = -8.892300000000006
$ws.Range("B57").Value = 5.163099999999995
$ws.Range("D57").Value = -8.333799999999997
$ws.Range("D61").Value = -7.7001
$ws.Range("D64").Value = -7.043299999999991
$ws.Range("C66").Value = -10.8751
$ws.Range("C67").Value = -10.49119999999999
$ws.Range("B69").Value = 5.404799999999994
$ws.Range("C73").Value = -10.67670000000001
$ws.Range("D75").Value = -8.454000000000006
$ws.Range("B76").Value = 4.713300000000001
$ws.Range("B78").Value = 10.54900000000001
$ws.Range("C78").Value = -10.91840000000001
$ws.Range("C80").Value = -10.93830000000001
$ws.Range("B82").Value = 5.1506
$ws.Range("D82").Value = -8.482300000000006
$ws.Range("B83").Value = 5.7028
$ws.Range("D90").Value = -7.206699999999993
$ws.Range("C91").Value = -12.97019999999999
$ws.Range("B93").Value = 4.675099999999996
$ws.Range("D93").Value = -7.212599999999993
$ws.Range("B97").Value = 5.411999999999995
$ws.Range("C97").Value = -10.5472
$ws.Range("C99").Value = -12.8723
$ws.Range("D102").Value = -7.640900000000005
$ws.Range("C104").Value = -11.93970000000001
$ws.Range("D105").Value = -7.597899999999998
